$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D, shifting existing Compensation/fiberAeff/fiberAlphadB
# columns (old D:H) one place to the right (new E:I).
$ws.Columns("D").Insert()

# New header cell for the inserted "Span (km)" column, matching the
# bold/centered/bordered look of the other header cells.
$ws.Range("D1").Value = "Span (km)"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").HorizontalAlignment = -4108
$ws.Range("D1").VerticalAlignment = -4160
$ws.Range("D1").Borders.LineStyle = 1

# Fill the new Span (km) column with 25 for every data row (2-25).
$ws.Range("D2:D25").Value = 25
